# Fix: Corrección de ortografía y fecha en el documento de requerimientos
$wb = $excel.ActiveWorkbook

$wsMatriz  = $wb.Worksheets.Item("Matriz de priorización")
$wsLeyenda = $wb.Worksheets.Item("Leyenda")

# --- Matriz de priorización sheet edits ---

# C1: "Última modificación: 30/enero/2025" -> "Fecha: 30/enero/2025"
$wsMatriz.Range("C1").Value = "Fecha: 30/enero/2025"

# C43: fix typo "terndrá" -> "tendrá"
$wsMatriz.Range("C43").Value = "El sistema tendrá una respuesta menor a dos segundos para cada página"

# C25: fix typo "recibira" -> "recibirá"
$wsMatriz.Range("C25").Value = "El cliente recibirá una alerta de los productos de la tienda que tengan un stock menor a cinco unidades"

# Selection / scroll state left on the Matriz sheet after editing
$wsMatriz.Range("G24").Select()
$wsMatriz.Application.ActiveWindow.ScrollRow = 20

# --- Leyenda sheet edits ---

# K5: fix typos "impresindible" -> "imprescindible", "acabo" -> "cabo"
$wsLeyenda.Range("K5").Value = "Requisito totalmente imprescindible que tiene que estar incluido ya que si no se llevan a cabo el proyecto no puede seguir adelante"

# K7: fix typo "inclurise" -> "incluirse"
$wsLeyenda.Range("K7").Value = "Requisitos que podrían incluirse si no afecta a nada más, es decir, son requisitos que sería bueno y podrían incluirse porque no cuesta demasiado implementarlos"

# Leave Leyenda as the active sheet/tab with L7 selected
$wsLeyenda.Activate()
$wsLeyenda.Range("L7").Select()
